# Adds two new columns (F, G) of commentary to the accuracy sheet:
#   F1 / G1 -> new header labels
#   F4 / G4 -> the single explanatory note that goes with row 4 (fish(3))
# Also widens the new columns and refreshes the view (zoom + selection)
# to match how the sheet was left after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns F & G
$ws.Range("F1").Value = "脊椎抓不好的原因"
$ws.Range("G1").Value = "如何改善"

# New explanatory note tied to row 4
$ws.Range("F4").Value = "沒有完整抓出脊椎"
$ws.Range("G4").Value = "用多階grayseperate將脊椎顯現出來"

# Widen the new columns so the text is readable
$ws.Columns("F").ColumnWidth = 17.426339285714285
$ws.Columns("G").ColumnWidth = 63.570870535714285

# Zoom out a bit and move the selection to the top of the new columns
$excel.ActiveWindow.Zoom = 115
[void]$ws.Range("G5").Select()
